# Apply updated TPM-based values to the Il6-Il6ra NATMI output sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.3056103333333333
$ws.Range("H2").Value = 0.916831
$ws.Range("I2").Value = 0.01726097181671177
$ws.Range("J2").Value = 0.01726097181671177
$ws.Range("M2").Value = 0.9317853333333334
$ws.Range("N2").Value = 2.795356
$ws.Range("O2").Value = 0.1255826100074751
$ws.Range("P2").Value = 0.1255826100074751
$ws.Range("Q2").Value = 0.2847632263151111
$ws.Range("R2").Value = 2.562869036836
$ws.Range("S2").Value = 0.002167677892008134
$ws.Range("T2").Value = 0.002167677892008134
# Row 3
$ws.Range("G3").Value = 0.3056103333333333
$ws.Range("H3").Value = 0.916831
$ws.Range("I3").Value = 0.01726097181671177
$ws.Range("J3").Value = 0.01726097181671177
$ws.Range("O3").Value = 0.7447810673036616
$ws.Range("P3").Value = 0.7447810673036616
$ws.Range("Q3").Value = 1.688818695607444
$ws.Range("R3").Value = 15.199368260467
$ws.Range("S3").Value = 0.01285564501234901
$ws.Range("T3").Value = 0.01285564501234901
# Row 4
$ws.Range("G4").Value = 0.3056103333333333
$ws.Range("H4").Value = 0.916831
$ws.Range("I4").Value = 0.01726097181671177
$ws.Range("J4").Value = 0.01726097181671177
$ws.Range("M4").Value = 0.9618626666666666
$ws.Range("O4").Value = 0.1296363226888633
$ws.Range("P4").Value = 0.1296363226888633
$ws.Range("Q4").Value = 0.2939551701808888
$ws.Range("R4").Value = 2.645596531628
$ws.Range("S4").Value = 0.002237648912354622
$ws.Range("T4").Value = 0.002237648912354622
# Row 5
$ws.Range("I5").Value = 0.8433360339088308
$ws.Range("J5").Value = 0.8433360339088307
$ws.Range("M5").Value = 0.9317853333333334
$ws.Range("N5").Value = 2.795356
$ws.Range("O5").Value = 0.1255826100074751
$ws.Range("P5").Value = 0.1255826100074751
$ws.Range("Q5").Value = 13.91295301526178
$ws.Range("R5").Value = 125.216577137356
$ws.Range("S5").Value = 0.1059083402516235
$ws.Range("T5").Value = 0.1059083402516235
# Row 6
$ws.Range("I6").Value = 0.8433360339088308
$ws.Range("J6").Value = 0.8433360339088307
$ws.Range("O6").Value = 0.7447810673036616
$ws.Range("P6").Value = 0.7447810673036616
$ws.Range("S6").Value = 0.6281007114302559
$ws.Range("T6").Value = 0.6281007114302558
# Row 7
$ws.Range("I7").Value = 0.8433360339088308
$ws.Range("J7").Value = 0.8433360339088307
$ws.Range("M7").Value = 0.9618626666666666
$ws.Range("O7").Value = 0.1296363226888633
$ws.Range("P7").Value = 0.1296363226888633
$ws.Range("S7").Value = 0.1093269822269513
$ws.Range("T7").Value = 0.1093269822269513
# Row 8
$ws.Range("I8").Value = 0.1394029942744574
$ws.Range("J8").Value = 0.1394029942744574
$ws.Range("M8").Value = 0.9317853333333334
$ws.Range("N8").Value = 2.795356
$ws.Range("O8").Value = 0.1255826100074751
$ws.Range("P8").Value = 0.1255826100074751
$ws.Range("Q8").Value = 2.299803674388
$ws.Range("R8").Value = 20.698233069492
$ws.Range("S8").Value = 0.01750659186384347
$ws.Range("T8").Value = 0.01750659186384347
# Row 9
$ws.Range("I9").Value = 0.1394029942744574
$ws.Range("J9").Value = 0.1394029942744574
$ws.Range("O9").Value = 0.7447810673036616
$ws.Range("P9").Value = 0.7447810673036616
$ws.Range("S9").Value = 0.1038247108610566
$ws.Range("T9").Value = 0.1038247108610566
# Row 10
$ws.Range("I10").Value = 0.1394029942744574
$ws.Range("J10").Value = 0.1394029942744574
$ws.Range("M10").Value = 0.9618626666666666
$ws.Range("O10").Value = 0.1296363226888633
$ws.Range("P10").Value = 0.1296363226888633
$ws.Range("S10").Value = 0.01807169154955732
$ws.Range("T10").Value = 0.01807169154955732
